$d = $word.ActiveDocument

# 1) Justify the first three body paragraphs (those that originally had no
#    explicit alignment) -- adds <w:jc w:val="both"/> to their pPr.
for ($i = 1; $i -le 3; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Alignment = 3   # wdAlignParagraphJustify
}

# 2) Remove the trailing "Repositorio git..." section: two empty
#    paragraphs, the "Repositorio git..." paragraph, the bold repo-name
#    paragraph, and the hyperlink paragraph. These are the five
#    paragraphs that follow the last image paragraph ("Listado de
#    componentes automatizados FE." image), right before the final
#    section break.
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "Repositorio") {
        $targetIdx = $i
        break
    }
}

if ($targetIdx -ge 1) {
    # Walk back two more paragraphs to include the two empty ones before it.
    $deleteStartIdx = $targetIdx - 2
    if ($deleteStartIdx -lt 1) { $deleteStartIdx = 1 }

    $rangeStart = $d.Paragraphs.Item($deleteStartIdx).Range.Start
    $rangeEnd = $d.Paragraphs.Item($d.Paragraphs.Count).Range.End

    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
